$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name to reflect the new "through" date
$ws.Name = "Through 2021-10-04"

# Update the row label for October to reflect the new "through" date
$ws.Range("A11").Value = "October (through 10-04)"

# Update October row (row 11) figures for the columns that changed
$ws.Range("B11").Value = 4
$ws.Range("E11").Value = 10
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 22
$ws.Range("H11").Value = 31

# Update Total row (row 12) figures for the columns that changed
$ws.Range("B12").Value = 230
$ws.Range("E12").Value = 558
$ws.Range("F12").Value = 425
$ws.Range("G12").Value = 923
$ws.Range("H12").Value = 1279
